# Login.xlsx - "Registor" sheet updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registor")

# Row 5: update the verification message (inserted "123") and flip runMode to "Y"
$ws.Range("J5").Value = "Welcome123 to your account. Here you can manage all of your personal information and orders."
$ws.Range("K5").Value = "Y"

# Move the sheet's saved selection from F2:F5 to the whole of row 3
$ws.Rows.Item(3).Select()
